$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the distinguishing data among rows 3, 4 and 6:
#   new row 3 <- old row 6
#   new row 4 <- old row 3
#   new row 6 <- old row 4
# Columns B, C, D, E, F, G, H, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AT, AW, AX, AY are identical across the three rows and stay untouched.

# Column I holds text-formatted numerals in the source file (t="inlineStr"),
# so force the "@" text number format before writing to stop the engine
# from coercing the numeric-looking strings into real numbers.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I6").NumberFormat = "@"

# Row 3 becomes what row 6 used to be
$ws.Range("A3").Value = 111416521
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "stjälkar/strån/skott"
$ws.Range("K3").Value = "blomning"
$ws.Range("Q3").Value = 359101.3469427949
$ws.Range("R3").Value = 6393205.997596246
$ws.Range("AC3").Value = "även ca 30 bladrosetter"

# Row 4 becomes what row 3 used to be
$ws.Range("A4").Value = 111416528
$ws.Range("I4").Value = "30"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("K4").ClearContents()
$ws.Range("Q4").Value = 359092.1819271583
$ws.Range("R4").Value = 6393204.710604292
$ws.Range("AC4").ClearContents()

# Row 6 becomes what row 4 used to be
$ws.Range("A6").Value = 111416525
$ws.Range("I6").Value = "4"
$ws.Range("J6").Value = "stjälkar/strån/skott"
$ws.Range("K6").Value = "blomning"
$ws.Range("Q6").Value = 359095.1406046218
$ws.Range("R6").Value = 6393212.639220579
$ws.Range("AC6").Value = "även bladrosetter på 1 kvm"
